$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45997
$ws.Range("B2").Value = 33.83
$ws.Range("C2").Value = 22.84
$ws.Range("D2").Value = 19.28
$ws.Range("E2").Value = 15.14
$ws.Range("F2").Value = 12.7
$ws.Range("G2").Value = 13.67
$ws.Range("H2").Value = 17.62
$ws.Range("I2").Value = 23.35
$ws.Range("J2").Value = 25.06
$ws.Range("K2").Value = 30.19
$ws.Range("L2").Value = 26.24
$ws.Range("M2").Value = 17.6
$ws.Range("N2").Value = 18.58
$ws.Range("O2").Value = 15.32
$ws.Range("P2").Value = 16.61
$ws.Range("Q2").Value = 19.51
$ws.Range("R2").Value = 30.37
$ws.Range("S2").Value = 52
$ws.Range("T2").Value = 61.05
$ws.Range("U2").Value = 64.95
$ws.Range("V2").Value = 69.90000000000001
$ws.Range("W2").Value = 68.65000000000001
$ws.Range("X2").Value = 66.95
$ws.Range("Y2").Value = 52.37
$ws.Range("Z2").Value = 33.07
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 64.47
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 69.28
$ws.Range("AF2").Value = 63
$ws.Range("AG2").Value = "1h-16h"
